# "adapted tests for SQLSource"
#
# Row 1 used to hold plain text header labels ("col1".."col4") stored as
# shared strings. The adapted test now expects row 1 to contain actual
# typed sample values (matching the data types found in the rest of the
# sheet), so the SQLSource reader can infer column types from the first
# row instead of treating it as a text header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: was the text "col1" -> becomes the integer 0 (keeps its existing
# integer number format, style index unchanged).
$ws.Cells.Item(1, 2).Value = 0

# C1: was the text "col2" -> becomes the plain number 34.4.
$ws.Cells.Item(1, 3).Value = 34.4

# D1: was the text "col3" -> becomes a blank cell formatted like the
# date column below it (copy the number format from D2, then clear the
# cell's contents so it stays empty).
$ws.Cells.Item(2, 4).Copy()
$ws.Cells.Item(1, 4).PasteSpecial(-4122)
$ws.Cells.Item(1, 4).Value = ""

# E1: was the text "col4" -> becomes the text "11:11:11" entered with a
# leading apostrophe (quote prefix) while using a time number format,
# so it is stored as a quote-prefixed, time-formatted text value.
$ws.Cells.Item(1, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(1, 5).Value = "'11:11:11"

# Selection moved from A7 to D2.
$ws.Range("D2").Select()
